# Applies corrected stock-count/value figures to the Companywise Stock
# Report: adjusts Qty (F) and Value (G) for a number of line items, the
# swap of two item rows for "Rasna" SKUs (394/395 and 404/405), and
# recomputed Sub Total / Grand Total (B) rows that roll up the changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F27").Value = 42
$ws.Range("G27").Value = 1506.12

$ws.Range("F30").Value = 131
$ws.Range("G30").Value = 3691.58

$ws.Range("F31").Value = 39
$ws.Range("G31").Value = 1038.96

$ws.Range("B34").Value = 56856.63

$ws.Range("F44").Value = 29
$ws.Range("G44").Value = 1024.28

$ws.Range("F57").Value = 48
$ws.Range("G57").Value = 1695.36

$ws.Range("F58").Value = 64
$ws.Range("G58").Value = 4987.52

$ws.Range("B66").Value = 197050.27

$ws.Range("F96").Value = 15
$ws.Range("G96").Value = 13754.1

$ws.Range("B97").Value = 15282.93

$ws.Range("F114").Value = 31
$ws.Range("G114").Value = 1448.94

$ws.Range("F121").Value = 10
$ws.Range("G121").Value = 3077.7

$ws.Range("B123").Value = 71215

$ws.Range("F141").Value = 45
$ws.Range("G141").Value = 2342.7

$ws.Range("B147").Value = 20881.96

$ws.Range("F217").Value = 28
$ws.Range("G217").Value = 2080.4

$ws.Range("B218").Value = 73477.87

$ws.Range("F278").Value = 31
$ws.Range("G278").Value = 4201.74

$ws.Range("B295").Value = 111663.29

$ws.Range("F324").Value = 18
$ws.Range("G324").Value = 3083.94

$ws.Range("F326").Value = 12
$ws.Range("G326").Value = 3846.24

$ws.Range("B328").Value = -15223.7

$ws.Range("F365").Value = 10
$ws.Range("G365").Value = 553.3

$ws.Range("F366").Value = 53
$ws.Range("G366").Value = 2932.49

$ws.Range("F370").Value = 198
$ws.Range("G370").Value = 32866.02

$ws.Range("B372").Value = 56287.3

$ws.Range("F380").Value = 36
$ws.Range("G380").Value = 1977.12

$ws.Range("B389").Value = 56378.14

$ws.Range("F404").Value = 7
$ws.Range("G404").Value = 2428.72

$ws.Range("B417").Value = 166168.4

$ws.Range("F430").Value = 222
$ws.Range("G430").Value = 10274.16

$ws.Range("B438").Value = 23822.22

$ws.Range("F454").Value = 68
$ws.Range("G454").Value = 19241.28

$ws.Range("F455").Value = 40
$ws.Range("G455").Value = 8889.200000000001

$ws.Range("F456").Value = 147
$ws.Range("G456").Value = 39670.89

$ws.Range("B458").Value = 94185.28999999999

$ws.Range("F478").Value = 10
$ws.Range("G478").Value = 2217.4

$ws.Range("F481").Value = 33
$ws.Range("G481").Value = 1797.84

$ws.Range("B482").Value = 2371.72

$ws.Range("B496").Value = 60025
$ws.Range("E496").Value = 37.22
$ws.Range("F496").Value = -98
$ws.Range("G496").Value = -3217.34

$ws.Range("B497").Value = 64833
$ws.Range("E497").Value = 34.9
$ws.Range("F497").Value = 88
$ws.Range("G497").Value = 2889.04

$ws.Range("B506").Value = 64830
$ws.Range("E506").Value = 34.9
$ws.Range("F506").Value = 84
$ws.Range("G506").Value = 2757.72

$ws.Range("B507").Value = 60022
$ws.Range("E507").Value = 37.22
$ws.Range("F507").Value = -113
$ws.Range("G507").Value = -3709.79

$ws.Range("F520").Value = 14
$ws.Range("G520").Value = 383.6

$ws.Range("F522").Value = 77
$ws.Range("G522").Value = 2051.28

$ws.Range("B525").Value = 117029.64

$ws.Range("F527").Value = 41
$ws.Range("G527").Value = 1357.51

$ws.Range("F531").Value = 214
$ws.Range("G531").Value = 7085.54

$ws.Range("B535").Value = 22476.48

$ws.Range("F558").Value = 174
$ws.Range("G558").Value = 21201.9

$ws.Range("B561").Value = 25186.26

$ws.Range("F563").Value = 17
$ws.Range("G563").Value = 3176.62

$ws.Range("F565").Value = 13
$ws.Range("G565").Value = 3652.35

$ws.Range("F569").Value = 1
$ws.Range("G569").Value = 584.72

$ws.Range("F571").Value = 1
$ws.Range("G571").Value = 561.65

$ws.Range("B573").Value = 17672.3

$ws.Range("F612").Value = 232
$ws.Range("G612").Value = 34895.12

$ws.Range("F617").Value = 9
$ws.Range("G617").Value = 433.08

$ws.Range("F620").Value = 354
$ws.Range("G620").Value = 27820.86

$ws.Range("B628").Value = 204555.42

$ws.Range("F642").Value = 3
$ws.Range("G642").Value = 3923.82

$ws.Range("F643").Value = 1
$ws.Range("G643").Value = 1435.41

$ws.Range("B657").Value = 72750.7

$ws.Range("F660").Value = 48
$ws.Range("G660").Value = 1427.52

$ws.Range("B668").Value = 11357.96

$ws.Range("F674").Value = 664
$ws.Range("G674").Value = 108305.04

$ws.Range("B680").Value = 109317.59

$ws.Range("F703").Value = 6
$ws.Range("G703").Value = 2911.62

$ws.Range("F706").Value = 116
$ws.Range("G706").Value = 4546.04

$ws.Range("B713").Value = 63078.36

$ws.Range("B718").Value = 2525352.15

$ws.Range("B719").Value = 2525352.15
